# Weekly update: insert a new week's record at row 121 for the
# "Feria Lagunitas de Puerto Montt - Cebollín" series, pushing the
# existing rows 121-161 down to 122-162 (the oldest record, previously
# row 121, becomes the new row 162).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 121:161 down to 122:162, leaving a blank row 121 behind
# (formatting of row 121, incl. the date style on column D, is carried
# along automatically by Excel's insert-shift behaviour).
$ws.Rows("121").Insert()

# Fill the new row 121 with this week's data. The descriptive columns
# (Mercado ID, Mercado, Región, Codreg, Categoría ID, Categoría,
# Variedad, Calidad, Unidad de comercialización, Origen, Kg o Unidades,
# Clasificación) are constant across every row of this subset.
$ws.Range("A121").Value = 4
$ws.Range("B121").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C121").Value = "Los Lagos"
$ws.Range("D121").Value = 44463
$ws.Range("E121").Value = 10
$ws.Range("F121").Value = 100112037
$ws.Range("G121").Value = "Cebollín"
$ws.Range("H121").Value = "Sin especificar"
$ws.Range("I121").Value = "Primera"
$ws.Range("J121").Value = 180
$ws.Range("K121").Value = 6000
$ws.Range("L121").Value = 6000
$ws.Range("M121").Value = 6000
$ws.Range("N121").Value = "$/paquete 36 unidades"
$ws.Range("O121").Value = "Región Metropolitana"
$ws.Range("P121").Value = 167
$ws.Range("Q121").Value = 36
$ws.Range("R121").Value = "Hortaliza"
